$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing US9 / US10 rows to their new, more descriptive names ---
$ws.Range("A12").Value = "US9PrintOrderToPDF"
$ws.Range("C12").Value = "PASS"

$ws.Range("A13").Value = "US9PrintOrderToPrinter"
$ws.Range("C13").Value = "PASS"

# --- Insert two new rows for the reprint test cases, right after row 13 ---
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = "US10RePrintOrderToPDF"
$ws.Range("C14").Value = "PASS"

$ws.Range("A15").Value = "US10RePrintOrderToPrinter"
$ws.Range("C15").Value = "PASS"

# --- Widen column A slightly to fit the longer test-case names ---
$ws.Columns.Item(1).ColumnWidth = 26

# --- Re-apply conditional formatting so it covers the two newly inserted rows ---
$oldCF = $ws.Range("C2:C629").FormatConditions
$oldCF.Item(1).ModifyAppliesToRange($ws.Range("C2:C631"))

# --- Selection, matching the saved cursor position from the edit ---
$ws.Range("C15").Select()
